$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.057.65'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.904.66'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.50'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5050'
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4092'
$ws.Range("E8").Value = '  +3.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08342'
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.19'
$ws.Range("E10").Value = '  -0.11%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.91'
$ws.Range("E12").Value = '  +5.59%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.380'
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.897.41'
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.207'
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.16'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001096'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06485'
$ws.Range("E19").Value = '  +2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.32'
$ws.Range("E20").Value = '  +2.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").Value = '  +2.21%  '
$ws.Range("D23").Value = '30.060.08'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '2.122.91'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.68'
$ws.Range("E27").Value = '  +3.61%  '
$ws.Range("E28").Value = '  +1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.282'
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.52'
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.148'
$ws.Range("E31").Value = '  +10.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1040'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.946'
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.747'
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02448'
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.356'
$ws.Range("E36").Value = '  +3.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06362'
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2136'
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6534'
$ws.Range("E39").Value = '  +4.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.189'
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.564'
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.37'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.213'
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.41'
$ws.Range("E44").Value = '  +3.75%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.203'
$ws.Range("E45").Value = '  +10.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6061'
$ws.Range("E46").Value = '  +2.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.612'
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.208'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.13'
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.154'
$ws.Range("E50").Value = '  +2.46%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.73'
$ws.Range("E51").Value = '  +2.09%  '
